# Fruta / hortaliza, semanal
# This edit re-orders the weekly price records (rows 2-43) by permuting the
# per-record columns: D (Fecha), I (Calidad), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg).
# Columns A,B,C,E,F,G,H,N,O,Q,R are identical for every record so they are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: target row -> source row (the target row receives the values that
# used to live in the source row before the edit).
$rowMap = @{
    2  = 9
    3  = 27
    4  = 31
    5  = 17
    6  = 16
    7  = 21
    8  = 33
    9  = 6
    10 = 10
    11 = 35
    12 = 11
    13 = 4
    14 = 32
    15 = 15
    16 = 19
    17 = 12
    18 = 28
    19 = 43
    20 = 36
    21 = 2
    22 = 38
    23 = 30
    24 = 24
    25 = 37
    26 = 20
    27 = 42
    28 = 41
    29 = 14
    30 = 18
    31 = 34
    32 = 8
    33 = 40
    34 = 29
    35 = 3
    36 = 39
    37 = 7
    38 = 23
    39 = 26
    40 = 5
    41 = 22
    42 = 13
    43 = 25
}

# Columns (by index) that move together with each record.
# D=4, I=9, J=10, K=11, L=12, M=13, P=16
$cols = @(4, 9, 10, 11, 12, 13, 16)

# First take a full snapshot of the original values so that reading and
# writing never interfere with each other (the permutation is not just a
# set of independent row shifts).
$snapshot = @{}
foreach ($r in 2..43) {
    foreach ($c in $cols) {
        $snapshot["$r-$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Now write back the permuted values.
foreach ($targetRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$targetRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value2 = $snapshot["$sourceRow-$c"]
    }
}
